$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9404744
$ws.Range("I32").Value = 8929810
$ws.Range("K32").Value = 8929810
$ws.Range("M32").Value = -8929523

$ws.Range("H45").Value = 810.44446
$ws.Range("I45").Value = 736.75
$ws.Range("J45").Value = 1400
$ws.Range("K45").Value = 736.75
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = -359.75
$ws.Range("N45").Value = -2154

$ws.Range("H74").Value = 4220.15
$ws.Range("I74").Value = 4682.2
$ws.Range("K74").Value = 4682.2
$ws.Range("M74").Value = -3808.2

$ws.Range("H77").Value = 4220.15
$ws.Range("I77").Value = 4682.2
$ws.Range("K77").Value = 23411
$ws.Range("M77").Value = -19043

$ws.Range("H110").Value = 4941.069
$ws.Range("I110").Value = 2285.9167
$ws.Range("K110").Value = 2285.9167
$ws.Range("M110").Value = -240.9167000000002

$ws.Range("H122").Value = 2663.9473
$ws.Range("I122").Value = 1975.9375
$ws.Range("J122").Value = 6333.3335
$ws.Range("K122").Value = 5927.8125
$ws.Range("L122").Value = 19000.0005
$ws.Range("M122").Value = -3477.8125
$ws.Range("N122").Value = -23900.0005

$ws.Range("H132").Value = 951740.2
$ws.Range("I132").Value = 1114647
$ws.Range("J132").Value = 164357.5
$ws.Range("K132").Value = 3343941
$ws.Range("L132").Value = 493072.5
$ws.Range("M132").Value = -3341411
$ws.Range("N132").Value = -498132.5

$ws.Range("H133").Value = 152779.4
$ws.Range("J133").Value = 152779.4
$ws.Range("L133").Value = 152779.4
$ws.Range("N133").Value = -157839.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2635310
$ws.Range("I107").Value = 3451356.8
$ws.Range("K107").Value = 3451356.8
$ws.Range("M107").Value = -3449436.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2099
$ws.Range("I10").Value = 2099
$ws.Range("K10").Value = 2099
$ws.Range("M10").Value = -1960

$ws.Range("H16").Value = 21743414
$ws.Range("I16").Value = 41669640
$ws.Range("K16").Value = 41669640
$ws.Range("M16").Value = -41669353

$ws.Range("H31").Value = 4580.3076
$ws.Range("I31").Value = 1278
$ws.Range("J31").Value = 7002
$ws.Range("K31").Value = 1278
$ws.Range("L31").Value = 7002
$ws.Range("M31").Value = -983
$ws.Range("N31").Value = -7592

$ws.Range("H34").Value = 4580.3076
$ws.Range("I34").Value = 1278
$ws.Range("J34").Value = 7002
$ws.Range("K34").Value = 1278
$ws.Range("L34").Value = 7002
$ws.Range("M34").Value = -1076
$ws.Range("N34").Value = -7406

$ws.Range("H105").Value = 76924480
$ws.Range("I105").Value = 76924480
$ws.Range("K105").Value = 76924480
$ws.Range("M105").Value = -76922733

$ws.Range("H107").Value = 1243.25
$ws.Range("I107").Value = 1022.6957
$ws.Range("K107").Value = 1022.6957
$ws.Range("M107").Value = 897.3043

$ws.Range("H113").Value = 21743414
$ws.Range("I113").Value = 41669640
$ws.Range("K113").Value = 41669640
$ws.Range("M113").Value = -41667470

$ws.Range("H122").Value = 2414.3667
$ws.Range("I122").Value = 2001.619
$ws.Range("J122").Value = 3377.4443
$ws.Range("K122").Value = 6004.857
$ws.Range("L122").Value = 10132.3329
$ws.Range("M122").Value = -3554.857
$ws.Range("N122").Value = -15032.3329

$ws.Range("H132").Value = 6833.0713
$ws.Range("I132").Value = 4335.892
$ws.Range("K132").Value = 13007.676
$ws.Range("M132").Value = -10477.676

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 921.3158
$ws.Range("I5").Value = 656.8182
$ws.Range("K5").Value = 1970.4546
$ws.Range("M5").Value = -1858.4546

$ws.Range("H105").Value = 28099.8
$ws.Range("J105").Value = 28099.8
$ws.Range("L105").Value = 84299.39999999999
$ws.Range("N105").Value = -89541.39999999999

$ws.Range("H135").Value = 921.3158
$ws.Range("I135").Value = 656.8182
$ws.Range("K135").Value = 5911.3638
$ws.Range("M135").Value = -3376.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1974.2916
$ws.Range("I97").Value = 1939.2778
$ws.Range("J97").Value = 2079.3333
$ws.Range("K97").Value = 1939.2778
$ws.Range("L97").Value = 2079.3333
$ws.Range("M97").Value = -1443.2778
$ws.Range("N97").Value = -3071.3333

$ws.Range("H112").Value = 75000
$ws.Range("J112").Value = 75000
$ws.Range("L112").Value = 75000
$ws.Range("N112").Value = -77216

$ws.Range("H113").Value = 6384.476
$ws.Range("I113").Value = 2223
$ws.Range("K113").Value = 2223
$ws.Range("M113").Value = -53

$ws.Range("H126").Value = 13164309
$ws.Range("I126").Value = 19234264
$ws.Range("J126").Value = 12741.25
$ws.Range("K126").Value = 57702792
$ws.Range("L126").Value = 38223.75
$ws.Range("M126").Value = -57700322
$ws.Range("N126").Value = -43163.75

$ws.Range("H132").Value = 35718756
$ws.Range("I132").Value = 45458708
$ws.Range("K132").Value = 136376124
$ws.Range("M132").Value = -136373594

$ws.Range("H136").Value = 16401.291
$ws.Range("J136").Value = 16401.291
$ws.Range("L136").Value = 49203.87300000001
$ws.Range("N136").Value = -54303.87300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 5566423
$ws.Range("I11").Value = 10004758
$ws.Range("J11").Value = 18504.25
$ws.Range("K11").Value = 10004758
$ws.Range("L11").Value = 18504.25
$ws.Range("M11").Value = -10004618
$ws.Range("N11").Value = -18784.25

$ws.Range("H61").Value = 5195.0586
$ws.Range("I61").Value = 3759.7917
$ws.Range("K61").Value = 3759.7917
$ws.Range("M61").Value = -3557.7917

$ws.Range("H100").Value = 2114.238
$ws.Range("I100").Value = 2100
$ws.Range("K100").Value = 2100
$ws.Range("M100").Value = -1559

$ws.Range("H113").Value = 5195.0586
$ws.Range("I113").Value = 3759.7917
$ws.Range("K113").Value = 3759.7917
$ws.Range("M113").Value = -1589.7917

$ws.Range("H122").Value = 5887.7
$ws.Range("I122").Value = 5334.3335
$ws.Range("K122").Value = 16003.0005
$ws.Range("M122").Value = -13553.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1597.9
$ws.Range("I81").Value = 1498.5
$ws.Range("J81").Value = 1829.8334
$ws.Range("K81").Value = 2997
$ws.Range("L81").Value = 3659.6668
$ws.Range("M81").Value = -1936
$ws.Range("N81").Value = -5781.6668

$ws.Range("H84").Value = 1597.9
$ws.Range("I84").Value = 1498.5
$ws.Range("J84").Value = 1829.8334
$ws.Range("K84").Value = 14985
$ws.Range("L84").Value = 18298.334
$ws.Range("M84").Value = -9681
$ws.Range("N84").Value = -28906.334

$ws.Range("H96").Value = 2288.9285
$ws.Range("I96").Value = 2125
$ws.Range("J96").Value = 2507.5
$ws.Range("K96").Value = 2125
$ws.Range("L96").Value = 2507.5
$ws.Range("M96").Value = -752
$ws.Range("N96").Value = -5253.5

$ws.Range("H107").Value = 6897328
$ws.Range("I107").Value = 9091529
$ws.Range("J107").Value = 1268.1428
$ws.Range("K107").Value = 27274587
$ws.Range("L107").Value = 3804.4284
$ws.Range("M107").Value = -27272667
$ws.Range("N107").Value = -7644.428400000001

$ws.Range("H132").Value = 7316.4614
$ws.Range("I132").Value = 6238.174
$ws.Range("K132").Value = 18714.522
$ws.Range("M132").Value = -16184.522
